# Fragen und Glossar ergaenzt
# Add a new glossary entry (FBM -> Fliessband-Modul) as a new row right
# after the existing last data row, matching the formatting (font/style
# and row height) of the preceding data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (A8:B8) down onto the
# new row before writing its values, so the new cells pick up the same
# cell style used by the other glossary rows.
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A9").Value = "FBM"
$ws.Range("B9").Value = "Fließband-Modul"

$ws.Rows.Item(9).RowHeight = 15.75

$ws.Range("B11").Select()
